# This script re-applies a set of full-row re-shuffles to the "Artfynd"
# worksheet. The underlying source data (sightings) were re-ordered: several
# rows exchanged their entire contents (all columns A:AY) with one another,
# while staying on the same row number. We therefore capture the full row
# values first (before any writes happen) and then redistribute them
# according to the mapping below, so that:
#
#   new row 46 = old row 47
#   new row 47 = old row 48
#   new row 48 = old row 46
#
#   new row 57 = old row 59
#   new row 59 = old row 57
#
#   new row 66 = old row 67
#   new row 67 = old row 66
#
#   new row 70 = old row 72
#   new row 71 = old row 70
#   new row 72 = old row 73
#   new row 73 = old row 71
#
#   new row 82 = old row 83
#   new row 83 = old row 82
#
#   new row 84 = old row 86
#   new row 86 = old row 84
#
#   new row 87 = old row 88
#   new row 88 = old row 87

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A .. AY hold data for each observation row.
$firstCol = "A"
$lastCol = "AY"

# A handful of columns store plain-text values that nevertheless look like
# numbers or ISO dates ("Antal" counts such as "1"/"2", and "Startdatum"/
# "Slutdatum" values such as "2026-02-04"). Assigning such a string straight
# to .Value2 makes Excel silently reinterpret it as a real number/date
# serial, which would corrupt the data. To keep them as genuine text we
# re-assign those specific columns using their .Text representation with a
# leading apostrophe (quote-prefix), which forces a literal text value.
$forceTextCols = @(9, 25, 27)   # I = Antal, Y = Startdatum, AA = Slutdatum

# Mapping: destination row -> source row (content that should end up there)
$rowMap = [ordered]@{
    46 = 47
    47 = 48
    48 = 46
    57 = 59
    59 = 57
    66 = 67
    67 = 66
    70 = 72
    71 = 70
    72 = 73
    73 = 71
    82 = 83
    83 = 82
    84 = 86
    86 = 84
    87 = 88
    88 = 87
}

# Capture the full row contents BEFORE making any changes, keyed by the
# original row number, so overlapping cycles do not clobber data we still
# need to read later.
$originalRows = @{}
$originalText = @{}
$uniqueSourceRows = $rowMap.Values | Sort-Object -Unique
foreach ($srcRow in $uniqueSourceRows) {
    $rng = $ws.Range("$firstCol$srcRow`:$lastCol$srcRow")
    $originalRows[$srcRow] = $rng.Value2

    $textVals = @{}
    foreach ($colIdx in $forceTextCols) {
        $textVals[$colIdx] = $ws.Cells.Item($srcRow, $colIdx).Text
    }
    $originalText[$srcRow] = $textVals
}

# Now write each destination row using the captured snapshot.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $destRng = $ws.Range("$firstCol$destRow`:$lastCol$destRow")
    $destRng.Value2 = $originalRows[$srcRow]

    foreach ($colIdx in $forceTextCols) {
        $txt = $originalText[$srcRow][$colIdx]
        $cell = $ws.Cells.Item($destRow, $colIdx)
        if ([string]::IsNullOrEmpty($txt)) {
            $cell.Value2 = $null
        } else {
            $cell.Value2 = "'" + $txt
        }
    }
}
